$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they stay stored
# as text (matching the source workbook) instead of being parsed as numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = '34.547.07'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.801.67'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '224.41'
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = '41.58'
$ws.Range("E8").Value = '  +14.72%  '
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("D11").Value = '0.0997'
$ws.Range("E11").Value = '  +3.45%  '
$ws.Range("D12").Value = '2.061.03'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '1.793.13'
$ws.Range("E13").Value = '  -0.16%  '
$ws.Range("D14").Value = '10.93'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("D15").Value = '34.466.96'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '67.24'
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("D19").Value = '240.24'
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").Value = '0.0₃0766'
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").Value = '11.14'
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").Value = '4.27'
$ws.Range("E23").Value = '  +4.94%  '
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("D25").Value = '171.91'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("E26").Value = '  -2.81%  '
$ws.Range("D27").Value = '17.37'
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '0.121'
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").Value = '3.87'
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").Value = '87.56'
$ws.Range("E35").Value = '  +8.02%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.647'
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.318.22'
$ws.Range("E37").Value = '  -3.08%  '
$ws.Range("D38").Value = '1.05'
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("E39").Value = '  +0.67%  '
$ws.Range("D40").Value = '14.75'
$ws.Range("E40").Value = '  +12.06%  '
$ws.Range("D41").Value = '2.34'
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("E42").Value = '  +5.08%  '
$ws.Range("D43").Value = '2.43'
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = '0.935'
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").Value = '0.0518'
$ws.Range("E46").Value = '  +4.26%  '
$ws.Range("D47").Value = '1.963.36'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '5.78'
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("D50").Value = '100.73'
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.0608'
$ws.Range("E51").Value = '  +0.58%  '
